# MR-B (EFGH) Attendance Sheet
# Clear attendance marks ("P"/"A") recorded for Session 19 (column Y) and
# Session 20 (column Z) for every student row (7 through 76). The
# Total-Absence (E) and Total-Present (F) formulas recompute automatically
# once the underlying attendance cells are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Y7:Z76").ClearContents()

# Re-apply the plain "blank" cell formatting that the rest of the blank
# cells in the sheet already use (same formatting as column AA), instead
# of leaving the old "has data" formatting behind.
$ws.Range("AA7").Copy()
$ws.Range("Y7:Z7").PasteSpecial(-4122)

$ws.Range("AA8").Copy()
$ws.Range("Y8:Z76").PasteSpecial(-4122)

$excel.CutCopyMode = 0
